$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 8907.5
$ws.Range("J26").Value = 10015
$ws.Range("L26").Value = 10015
$ws.Range("N26").Value = -10703
$ws.Range("H69").Value = 3862.2222
$ws.Range("J69").Value = 3293.3333
$ws.Range("L69").Value = 9879.999899999999
$ws.Range("N69").Value = -11627.9999
$ws.Range("H72").Value = 3862.2222
$ws.Range("J72").Value = 3293.3333
$ws.Range("L72").Value = 29639.9997
$ws.Range("N72").Value = -38375.9997
$ws.Range("H100").Value = 1589
$ws.Range("I100").Value = 760.8
$ws.Range("J100").Value = 2417.2
$ws.Range("K100").Value = 760.8
$ws.Range("L100").Value = 2417.2
$ws.Range("M100").Value = -219.8
$ws.Range("N100").Value = -3499.2
$ws.Range("H106").Value = 3879
$ws.Range("I106").Value = 3899.8333
$ws.Range("J106").Value = 3854
$ws.Range("K106").Value = 3899.8333
$ws.Range("L106").Value = 3854
$ws.Range("M106").Value = -3268.8333
$ws.Range("N106").Value = -5116
$ws.Range("H138").Value = 2555.6575
$ws.Range("I138").Value = 1311.1562
$ws.Range("J138").Value = 3526.9756
$ws.Range("K138").Value = 3933.4686
$ws.Range("L138").Value = 10580.9268
$ws.Range("M138").Value = 1206.5314
$ws.Range("N138").Value = -20860.9268

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1685.591
$ws.Range("I61").Value = 1350.8823
$ws.Range("J61").Value = 2823.6
$ws.Range("K61").Value = 1350.8823
$ws.Range("L61").Value = 2823.6
$ws.Range("M61").Value = -1138.8823
$ws.Range("N61").Value = -3247.6
$ws.Range("H74").Value = 846.1739
$ws.Range("I74").Value = 737.44446
$ws.Range("J74").Value = 1237.6
$ws.Range("K74").Value = 737.44446
$ws.Range("L74").Value = 1237.6
$ws.Range("M74").Value = 136.55554
$ws.Range("N74").Value = -2985.6
$ws.Range("H77").Value = 846.1739
$ws.Range("I77").Value = 737.44446
$ws.Range("J77").Value = 1237.6
$ws.Range("K77").Value = 3687.2223
$ws.Range("L77").Value = 6188
$ws.Range("M77").Value = 680.7776999999996
$ws.Range("N77").Value = -14924
$ws.Range("H92").Value = 33000
$ws.Range("J92").Value = 33000
$ws.Range("L92").Value = 33000
$ws.Range("N92").Value = -37992
$ws.Range("H136").Value = 1685.591
$ws.Range("I136").Value = 1350.8823
$ws.Range("J136").Value = 2823.6
$ws.Range("K136").Value = 4052.6469
$ws.Range("L136").Value = 8470.799999999999
$ws.Range("M136").Value = -1502.6469
$ws.Range("N136").Value = -13570.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66786.586
$ws.Range("I86").Value = 80342.28999999999
$ws.Range("J86").Value = 3526.6667
$ws.Range("K86").Value = 80342.28999999999
$ws.Range("L86").Value = 3526.6667
$ws.Range("M86").Value = -79219.28999999999
$ws.Range("N86").Value = -5772.6667
$ws.Range("H89").Value = 66786.586
$ws.Range("I89").Value = 80342.28999999999
$ws.Range("J89").Value = 3526.6667
$ws.Range("K89").Value = 401711.45
$ws.Range("L89").Value = 17633.3335
$ws.Range("M89").Value = -396095.45
$ws.Range("N89").Value = -28865.3335
$ws.Range("H134").Value = 10389.051
$ws.Range("I134").Value = 10695.607
$ws.Range("K134").Value = 32086.821
$ws.Range("M134").Value = -29551.821

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35643.188
$ws.Range("I31").Value = 757.61536
$ws.Range("J31").Value = 50760.266
$ws.Range("K31").Value = 757.61536
$ws.Range("L31").Value = 50760.266
$ws.Range("M31").Value = -462.61536
$ws.Range("N31").Value = -51350.266
$ws.Range("H34").Value = 35643.188
$ws.Range("I34").Value = 757.61536
$ws.Range("J34").Value = 50760.266
$ws.Range("K34").Value = 757.61536
$ws.Range("L34").Value = 50760.266
$ws.Range("M34").Value = -555.61536
$ws.Range("N34").Value = -51164.266
$ws.Range("H47").Value = 18640
$ws.Range("I47").Value = 16700
$ws.Range("J47").Value = 21550
$ws.Range("K47").Value = 16700
$ws.Range("L47").Value = 21550
$ws.Range("M47").Value = -16134
$ws.Range("N47").Value = -22682

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 782.7692
$ws.Range("I2").Value = 15.5
$ws.Range("J2").Value = 1123.7778
$ws.Range("K2").Value = 93
$ws.Range("L2").Value = 6742.666800000001
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = -6968.666800000001
$ws.Range("H113").Value = 942.72
$ws.Range("J113").Value = 573.3333
$ws.Range("L113").Value = 1719.9999
$ws.Range("N113").Value = -6059.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 5509
$ws.Range("J40").Value = 5509
$ws.Range("L40").Value = 5509
$ws.Range("N40").Value = -5811

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1722.5714
$ws.Range("J22").Value = 285.33334
$ws.Range("L22").Value = 285.33334
$ws.Range("N22").Value = -875.33334
$ws.Range("H27").Value = 1722.5714
$ws.Range("J27").Value = 285.33334
$ws.Range("L27").Value = 285.33334
$ws.Range("N27").Value = -499.33334
$ws.Range("H82").Value = 2290.3333
$ws.Range("H85").Value = 2290.3333
$ws.Range("H136").Value = 1918.9375
$ws.Range("I136").Value = 1530.4
$ws.Range("K136").Value = 4591.200000000001
$ws.Range("M136").Value = -2041.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 18377.5
$ws.Range("J40").Value = 18377.5
$ws.Range("L40").Value = 18377.5
$ws.Range("N40").Value = -18675.5
$ws.Range("H81").Value = 182532.73
$ws.Range("I81").Value = 167213.5
$ws.Range("J81").Value = 200915.8
$ws.Range("K81").Value = 334427
$ws.Range("L81").Value = 401831.6
$ws.Range("M81").Value = -333366
$ws.Range("N81").Value = -403953.6
$ws.Range("H84").Value = 182532.73
$ws.Range("I84").Value = 167213.5
$ws.Range("J84").Value = 200915.8
$ws.Range("K84").Value = 1672135
$ws.Range("L84").Value = 2009158
$ws.Range("M84").Value = -1666831
$ws.Range("N84").Value = -2019766
$ws.Range("H126").Value = 1365.6957
$ws.Range("I126").Value = 1414.2222
$ws.Range("J126").Value = 1191
$ws.Range("K126").Value = 4242.6666
$ws.Range("L126").Value = 3573
$ws.Range("M126").Value = -1772.6666
$ws.Range("N126").Value = -8513
$ws.Range("H136").Value = 21025.596
$ws.Range("I136").Value = 50652.7
$ws.Range("J136").Value = 5010.946
$ws.Range("K136").Value = 151958.1
$ws.Range("L136").Value = 15032.838
$ws.Range("M136").Value = -149408.1
$ws.Range("N136").Value = -20132.838
